$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.641.55'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').Value = '2.251.89'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.632'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '76.60'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.01%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.629'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '44.76'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.94%  '
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.38'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.70%  '
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.77'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.868'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').Value = '2.270.39'
$ws.Range('E16').Value = '  +0.57%  '
$ws.Range('D17').Value = '42.470.45'
$ws.Range('E17').Value = '  +1.24%  '
$ws.Range('E18').Value = '  +3.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +56.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.51%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.85%  '
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.76'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.02%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.69'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0826'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +18.80%  '
$ws.Range('E34').Value = '  -0.94%  '
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.77'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0319'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '14.31'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.23%  '
$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '64.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.26%  '
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '108.48'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.94'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.29%  '
$ws.Range('E45').Value = '  +2.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.20'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.02%  '
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.425'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.12%  '
